$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark that used to sit in the MERGE
#    section (right after the word "branch" in "Mesclagem entre
#    branch"). Pre-existing bookmarks can be located by name but can't
#    be mutated/deleted directly through this object model, so instead
#    we delete the whole paragraph's text and retype it identically --
#    that drops the now-orphaned bookmark markers along with it.
#    This must run before step 2 below: once a *new* bookmark named
#    "_GoBack" exists, name-based lookups resolve to it first (it is
#    earlier in document order), so the stale one becomes unreachable.
# ---------------------------------------------------------------------
$goBackOld = $d.Bookmarks("_GoBack")
$mergePara = $goBackOld.Range.Paragraphs(1).Range
$mergeText = $mergePara.Text
$mergeStart = $mergePara.Start
$mergePara.Delete()
$d.Range($mergeStart, $mergeStart).InsertBefore($mergeText)

# ---------------------------------------------------------------------
# 2) Add a new "_GoBack" bookmark collapsed right before the very first
#    run of the document (before "TRABALHO SOBRE GIT").
#    Inserting a bookmark exactly at absolute position 0 makes the
#    runtime anchor bookmarkEnd in the *next* paragraph, so we work
#    around that by temporarily typing a placeholder character at the
#    start, anchoring the bookmark right after it (position 1), and
#    then deleting the placeholder again.
# ---------------------------------------------------------------------
$startRange = $d.Range(0, 0)
$startRange.InsertBefore("X")
$afterPlaceholder = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $afterPlaceholder)
$d.Range(0, 1).Delete()

# ---------------------------------------------------------------------
# 3) Merge the four runs that make up the
#    " clean –f <caminho> (remove arquivo não rastreados de lugar
#    determinado)" sentence into a single run/text node.
#    A direct Range.Text assignment is a no-op when the replacement
#    text equals the existing text, so we first swap in a throwaway
#    marker (forcing the runs to collapse into one) and then replace
#    that marker with the real text.
# ---------------------------------------------------------------------
$cleanTarget = " clean –f <caminho> (remove arquivo não rastreados de lugar determinado)"
$cleanRange = $d.Content
$cleanRange.Find.Execute($cleanTarget, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$cleanRange.Text = "TEMP_MARKER_CLEAN_F"
$cleanRange2 = $d.Content
$cleanRange2.Find.Execute("TEMP_MARKER_CLEAN_F", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$cleanRange2.Text = $cleanTarget

# ---------------------------------------------------------------------
# 4) Merge the two runs " " + "apelido " (right after "git pull") into
#    a single run " apelido ".
# ---------------------------------------------------------------------
$pullRange = $d.Content
$pullRange.Find.Execute("pull", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$apelidoRange = $d.Range($pullRange.End, $pullRange.End + 9)
$apelidoRange.Text = "TEMP_MARKER_APELIDO"
$apelidoRange2 = $d.Content
$apelidoRange2.Find.Execute("TEMP_MARKER_APELIDO", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$apelidoRange2.Text = " apelido "

Write-Output "done"
